# Applies the "Editing Choice question, FormInstance generation and view,
# new window for filling" update to the WorkReport workbook: appends two
# new log rows (33 and 34) below the existing data on Sheet1, and extends
# the running-total SUM formula in C3 to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared strings are allocated in first-use order; the commit's new text
# ("Editing Choice question...") lands on index 28, and the shorter text
# ("Small improvements...") on index 29 - even though it's row 33's B cell
# that ends up pointing at index 29 and row 34's B cell at index 28. Touch
# B34 first so the string table gets built in that order.
$ws.Range("B34").Value = "Editing Choice question, FormInstance generation and view, new window for filling"
$ws.Range("B33").Value = "Small improvements (cropping form name, …)"

# New row 33 - copy the date formatting from the last existing row first so
# the new cell reuses the existing date style instead of minting a new one,
# then overwrite with the real value.
$ws.Range("C33").Value = 0.5
$ws.Range("D32").Copy($ws.Range("D33"))
$ws.Range("D33").Value = 40884

# New row 34
$ws.Range("C34").Value = 6
$ws.Range("D32").Copy($ws.Range("D34"))
$ws.Range("D34").Value = 40885

# Extend the total formula to include the two new rows
$ws.Range("C3").Formula = "=SUM(C4:C565)"

# Move the active selection to follow the newly-added last row
$ws.Range("B34").Select()
